$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update timestamp string
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 01:35"

# Row 4: Estados Unidos - updated counts
$ws.Range("B4").Value = 1429084
$ws.Range("C4").Value = 20448
$ws.Range("D4").Value = 307895
$ws.Range("E4").Value = 1036037
$ws.Range("G4").Value = 1727
$ws.Range("H4").Value = 85152

# Row 9: Brasil - updated counts
$ws.Range("B9").Value = 189157
$ws.Range("C9").Value = 11555
$ws.Range("E9").Value = 97584

# Row 11: Alemania - updated counts
$ws.Range("E11").Value = 17537
$ws.Range("G11").Value = 123
$ws.Range("H11").Value = 7861

# Rows 54-55: Argentina overtakes Malasia in ranking
$ws.Range("A54").Value = "Argentina"
$ws.Range("B54").Value = 6879
$ws.Range("C54").Value = 316
$ws.Range("D54").Value = 2266
$ws.Range("E54").Value = 4284
$ws.Range("F54").Value = 170
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 329

$ws.Range("A55").Value = "Malasia"
$ws.Range("B55").Value = 6779
$ws.Range("C55").Value = 37
$ws.Range("D55").Value = 5281
$ws.Range("E55").Value = 1387
$ws.Range("F55").Value = 16
$ws.Range("G55").Value = 2
$ws.Range("H55").Value = 111

# Rows 83-87: Sudan jumps ahead of Cuba, Islandia, Estonia, Republica de Macedonia
$ws.Range("A83").Value = "Sudan"
$ws.Range("B83").Value = 1818
$ws.Range("C83").Value = 157
$ws.Range("D83").Value = 198
$ws.Range("E83").Value = 1530
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 10
$ws.Range("H83").Value = 90

$ws.Range("A84").Value = "Cuba"
$ws.Range("B84").Value = 1810
$ws.Range("C84").Value = 6
$ws.Range("D84").Value = 1326
$ws.Range("E84").Value = 405
$ws.Range("F84").Value = 7
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 79

$ws.Range("A85").Value = "Islandia"
$ws.Range("B85").Value = 1802
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 1780
$ws.Range("E85").Value = 12
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 10

$ws.Range("A86").Value = "Estonia"
$ws.Range("B86").Value = 1751
$ws.Range("C86").Value = 5
$ws.Range("D86").Value = 777
$ws.Range("E86").Value = 913
$ws.Range("F86").Value = 5
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 61

$ws.Range("A87").Value = "Republica de Macedonia"
$ws.Range("B87").Value = 1694
$ws.Range("C87").Value = 20
$ws.Range("D87").Value = 1229
$ws.Range("E87").Value = 370
$ws.Range("F87").Value = 21
$ws.Range("G87").Value = 3
$ws.Range("H87").Value = 95

# Row 117: Uruguay - updated counts
$ws.Range("B117").Value = 719
$ws.Range("C117").Value = 2
$ws.Range("D117").Value = 545
$ws.Range("E117").Value = 155

# Row 162: Guyana - updated counts
$ws.Range("D162").Value = 41
$ws.Range("E162").Value = 62
$ws.Range("F162").Value = 3

